$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.007921659358885036, 104.7277239638582, -76.12851921821058, 86.83826899324438, 1.087736738464413, 0.137925699246686, 0.1321353430312686, -0.457978408653416, 4.353809811669684, 257, -62.91272041886447, 0, 4.355622319099083, 5.02228898576575),
    @(-8.0941154069587, 9.110530214514723, 1.378088573055627, 81.22453385091053, -0.6894683831849064, -0.02385453739126353, 0.5090568506303041, -0.7261686376390966, 4.35380761793208, 937, -18.4904538046267, 0, 4.355684915537104, 5.022351582203771),
    @(30.15554940028988, -159.0537175503449, 0.1741289810405534, 203.5577564041317, 0.05047289782313857, -0.2456751588185071, 0.7346961534988825, -0.3938585943787365, 4.353778286468835, 51, -18.48982860918638, 0, 4.355743674148179, 5.022410340814846),
    @(10.45465227168621, 6.169765644211846, [double]"4.689390985214435e-05", 75.44439061269631, 0.2488365208279393, -0.9673196855943211, 1.913238163245572, -0.5302732929203884, 4.353754389013545, 621, -33.27196597344029, 0, 4.355746050719752, 5.022412717386419),
    @(26.4088679892641, 92.41251403425821, 53.52965755472344, 0.00128606673047137, -1.538403256623806, -0.3838261788064343, 0.1135073258319177, 1.390909124276527, 4.353795121719429, 945, -100.2990946094875, 0, 4.35574622047838, 5.022412887145047),
    @(0.04551887775591026, 3.049401829752029, 52.39268169665311, 26.30273150456254, 0.6794030986343471, 0.3922754691060586, -0.7813892694683786, -0.4080308380812336, 4.353824055077897, 778, -18.48995982126471, 0, 4.355755810203156, 5.022422476869823),
    @(2.641637624666615, 1.266984588286954, 17.81970626781357, 65.82087572925903, 0.2133985500628559, 0.4984203616193201, -0.3325013672918351, -0.8117430424133312, 4.353805514627552, 469, -18.48968644373167, 0, 4.355859835452137, 5.022526502118803),
    @(25.91285118940834, 0.004368737047227582, 1.964315188081088, 81.24557095895548, -0.08647006918719491, 1.041943234144266, 0.4566069549973326, -0.9033456813570244, 4.353778999871952, 802, -28.72136330555168, 0, 4.355939097714391, 5.022605764381058),
    @(0.1779081014491591, -93.40243263581021, 170.7072061467497, 34.91794191928116, 0.7557727624322843, -1.04992582854416, -0.8621751461306277, 0.06158433865995638, 4.35371795186191, 526, -48.06299196929407, 0, 4.355946107051409, 5.022612773718076),
    @(1.348284501574201, [double]"2.714814759204617e-06", 92.84530405557746, 24.44006165666981, 0.5223596406124948, 1.697505301138365, -1.114819023889696, -0.2036030782800891, 4.353815200091578, 98, -18.50677047089281, 0, 4.356042197637365, 5.022708864304032)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
